$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append 12 new rows (64-75) of channel-label lookups below the existing
# module-coefficient table. Column F is intentionally left blank (it mirrors
# the existing table layout, where F holds a value only sporadically).
# ---------------------------------------------------------------------------

# Row 64 - "_1_v" aliases
$ws.Range("A64").Value = "str_1_v"
$ws.Range("B64").Value = "LM1_V_Avg"
$ws.Range("C64").Value = "PH1_V_Avg"
$ws.Range("D64").Value = "CP1_V_Avg"
$ws.Range("E64").Value = "CM1_V_Avg"
$ws.Range("G64").Value = "J1_V_Avg"
$ws.Range("H64").Value = "MM1_V_Avg"
$ws.Range("I64").Value = "HP1_V_Avg"
$ws.Range("J64").Value = "HM1_V_Avg"

# Row 65 - "_2_v" aliases
$ws.Range("A65").Value = "str_2_v"
$ws.Range("B65").Value = "LM2_V_Avg"
$ws.Range("C65").Value = "PH2_V_Avg"
$ws.Range("D65").Value = "CP2_V_Avg"
$ws.Range("E65").Value = "CM2_V_Avg"
$ws.Range("G65").Value = "J2_V_Avg"
$ws.Range("H65").Value = "MM2_V_Avg"
$ws.Range("I65").Value = "HP2_V_Avg"
$ws.Range("J65").Value = "HM2_V_Avg"

# Row 66 - "_3_v" aliases
$ws.Range("A66").Value = "str_3_v"
$ws.Range("B66").Value = "LM3_V_Avg"
$ws.Range("C66").Value = "PH3_V_Avg"
$ws.Range("D66").Value = "CP3_V_Avg"
$ws.Range("E66").Value = "CM3_V_Avg"
$ws.Range("G66").Value = "J3_V_Avg"
$ws.Range("H66").Value = "MM3_V_Avg"
$ws.Range("I66").Value = "HP3_V_Avg"
$ws.Range("J66").Value = "HM3_V_Avg"

# Row 67 - "_4_v" aliases
$ws.Range("A67").Value = "str_4_v"
$ws.Range("B67").Value = "LM4_V_Avg"
$ws.Range("C67").Value = "PH4_V_Avg"
$ws.Range("D67").Value = "CP4_V_Avg"
$ws.Range("E67").Value = "CM4_V_Avg"
$ws.Range("G67").Value = "J4_V_Avg"
$ws.Range("H67").Value = "MM4_V_Avg"
$ws.Range("I67").Value = "HP4_V_Avg"
$ws.Range("J67").Value = "HM4_V_Avg"

# Row 68 - "_1_i" aliases
$ws.Range("A68").Value = "str_1_i"
$ws.Range("B68").Value = "LM1_I_Avg"
$ws.Range("C68").Value = "PH1_I_Avg"
$ws.Range("D68").Value = "CP1_I_Avg"
$ws.Range("E68").Value = "CM1_I_Avg"
$ws.Range("G68").Value = "J1_I_Avg"
$ws.Range("H68").Value = "MM1_I_Avg"
$ws.Range("I68").Value = "HP1_I_Avg"
$ws.Range("J68").Value = "HM1_I_Avg"

# Row 69 - "_2_i" aliases
$ws.Range("A69").Value = "str_2_i"
$ws.Range("B69").Value = "LM2_I_Avg"
$ws.Range("C69").Value = "PH2_I_Avg"
$ws.Range("D69").Value = "CP2_I_Avg"
$ws.Range("E69").Value = "CM2_I_Avg"
$ws.Range("G69").Value = "J2_I_Avg"
$ws.Range("H69").Value = "MM2_I_Avg"
$ws.Range("I69").Value = "HP2_I_Avg"
$ws.Range("J69").Value = "HM2_I_Avg"

# Row 70 - "_3_i" aliases
$ws.Range("A70").Value = "str_3_i"
$ws.Range("B70").Value = "LM3_I_Avg"
$ws.Range("C70").Value = "PH3_I_Avg"
$ws.Range("D70").Value = "CP3_I_Avg"
$ws.Range("E70").Value = "CM3_I_Avg"
$ws.Range("G70").Value = "J3_I_Avg"
$ws.Range("H70").Value = "MM3_I_Avg"
$ws.Range("I70").Value = "HP3_I_Avg"
$ws.Range("J70").Value = "HM3_I_Avg"

# Row 71 - "_4_i" aliases
$ws.Range("A71").Value = "str_4_i"
$ws.Range("B71").Value = "LM4_I_Avg"
$ws.Range("C71").Value = "PH4_I_Avg"
$ws.Range("D71").Value = "CP4_I_Avg"
$ws.Range("E71").Value = "CM4_I_Avg"
$ws.Range("G71").Value = "J4_I_Avg"
$ws.Range("H71").Value = "MM4_I_Avg"
$ws.Range("I71").Value = "HP4_I_Avg"
$ws.Range("J71").Value = "HM4_I_Avg"

# Row 72 - "_1_rtd" aliases
$ws.Range("A72").Value = "str_1_rtd"
$ws.Range("B72").Value = "LM1_RTD_Avg"
$ws.Range("C72").Value = "PH1_RTD_Avg"
$ws.Range("D72").Value = "CP1_RTD_Avg"
$ws.Range("E72").Value = "CM1_RTD_Avg"
$ws.Range("G72").Value = "J1_RTD_Avg"
$ws.Range("H72").Value = "MM1_RTD_Avg"
$ws.Range("I72").Value = "HP1_RTD_Avg"
$ws.Range("J72").Value = "HM1_RTD_Avg"

# Row 73 - "_2_rtd" aliases
$ws.Range("A73").Value = "str_2_rtd"
$ws.Range("B73").Value = "LM2_RTD_Avg"
$ws.Range("C73").Value = "PH2_RTD_Avg"
$ws.Range("D73").Value = "CP2_RTD_Avg"
$ws.Range("E73").Value = "CM2_RTD_Avg"
$ws.Range("G73").Value = "J2_RTD_Avg"
$ws.Range("H73").Value = "MM2_RTD_Avg"
$ws.Range("I73").Value = "HP2_RTD_Avg"
$ws.Range("J73").Value = "HM2_RTD_Avg"

# Row 74 - "_3_rtd" aliases
$ws.Range("A74").Value = "str_3_rtd"
$ws.Range("B74").Value = "LM3_RTD_Avg"
$ws.Range("C74").Value = "PH3_RTD_Avg"
$ws.Range("D74").Value = "CP3_RTD_Avg"
$ws.Range("E74").Value = "CM3_RTD_Avg"
$ws.Range("G74").Value = "J3_RTD_Avg"
$ws.Range("H74").Value = "MM3_RTD_Avg"
$ws.Range("I74").Value = "HP3_RTD_Avg"
$ws.Range("J74").Value = "HM3_RTD_Avg"

# Row 75 - "_4_rtd" aliases
$ws.Range("A75").Value = "str_4_rtd"
$ws.Range("B75").Value = "LM4_RTD_Avg"
$ws.Range("C75").Value = "PH4_RTD_Avg"
$ws.Range("D75").Value = "CP4_RTD_Avg"
$ws.Range("E75").Value = "CM4_RTD_Avg"
$ws.Range("G75").Value = "J4_RTD_Avg"
$ws.Range("H75").Value = "MM4_RTD_Avg"
$ws.Range("I75").Value = "HP4_RTD_Avg"
$ws.Range("J75").Value = "HM4_RTD_Avg"

# ---------------------------------------------------------------------------
# Formatting for the new block: column A keeps the same bold style used by
# the rest of the table's label column (A1:A63); columns B:J get a new
# right-aligned style (also stamps the otherwise-empty F64 cell, matching
# the single stray styled-but-blank cell left behind by the original edit).
# ---------------------------------------------------------------------------
$ws.Range("A64:A75").Font.Bold = $true
$ws.Range("B64:J64").HorizontalAlignment = -4152
$ws.Range("B65:J75").HorizontalAlignment = -4152

# ---------------------------------------------------------------------------
# Restore the view: the sheet had scrolled down and the active selection had
# moved to L73 in the saved workbook.
# ---------------------------------------------------------------------------
$ws.Range("L73").Select() | Out-Null
